# Update cryptocurrency price and volume data (Price column D, Volume(1h) column E)
# Values are text in the source sheet; force text format for numeric-looking
# price strings so Excel does not auto-convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.271.60"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.270.34"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.95"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.11"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.70"
$ws.Range("E10").Value = "  +9.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.66"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "2.622.01"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.40"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "2.265.17"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  +2.46%  "
$ws.Range("D18").Value = "42.172.64"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.61"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.74"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.82"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.60"
$ws.Range("E28").Value = "  +6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.74"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.27"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +4.80%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.13"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.115"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.44"
$ws.Range("E42").Value = "  +14.31%  "
$ws.Range("D43").Value = "1.995.41"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0287"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.95"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.98"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.22"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.16"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.37"
$ws.Range("E51").Value = "  +0.09%  "
